$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.242.98"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "3.260.52"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'397.95"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'108.95"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").Value = "'0.579"
$ws.Range("E7").Value = "  +4.77%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D9").Value = "'0.619"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "'39.22"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.0948"
$ws.Range("E11").Value = "  +4.20%  "
$ws.Range("D12").Value = "'0.142"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "3.775.18"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "'8.25"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "3.250.50"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "'11.08"
$ws.Range("E18").Value = "  +3.38%  "
$ws.Range("D19").Value = "57.031.20"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  +4.84%  "
$ws.Range("D22").Value = "'12.89"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "'298.19"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").Value = "'74.00"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "'28.11"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").Value = "'7.85"
$ws.Range("E28").Value = "  -4.25%  "
$ws.Range("D29").Value = "'7.44"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "'11.17"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "'39.91"
$ws.Range("E34").Value = "  +9.75%  "
$ws.Range("D35").Value = "'0.0494"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").Value = "'51.37"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'3.46"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").Value = "'137.79"
$ws.Range("E41").Value = "  +3.78%  "
$ws.Range("D43").Value = "'0.286"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("D46").Value = "'16.68"
$ws.Range("E46").Value = "  -3.66%  "
$ws.Range("D47").Value = "'22.28"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("D49").Value = "2.144.72"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'2.46"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'1.97"
$ws.Range("E51").Value = "  -7.27%  "
